$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.395.51"
$ws.Range("E2").Value = "  +0.02%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.572.31"
$ws.Range("E3").Value = "  +0.28%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("E5").Value = "  +0.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.45"
$ws.Range("E6").Value = "  +0.38%  "

$ws.Range("E7").Value = "  +2.36%  "

$ws.Range("E8").Value = "  +1.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3423"
$ws.Range("E9").Value = "  +0.98%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07645"
$ws.Range("E10").Value = "  +0.38%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.152"
$ws.Range("E11").Value = "  -1.64%  "

$ws.Range("E12").Value = "  +0.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.18"
$ws.Range("E13").Value = "  -0.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.007"
$ws.Range("E14").Value = "  -0.84%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.951"
$ws.Range("E15").Value = "  +0.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.571.00"
$ws.Range("E16").Value = "  -0.08%  "

$ws.Range("E17").Value = "  -0.29%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.94"
$ws.Range("E18").Value = "  +0.86%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06741"
$ws.Range("E19").Value = "  -0.42%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  +0.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.78"
$ws.Range("E21").Value = "  +1.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.213"
$ws.Range("E22").Value = "  -0.56%  "

$ws.Range("E23").Value = "  -0.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.378.58"
$ws.Range("E24").Value = "  -0.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.398"
$ws.Range("E25").Value = "  +0.51%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.674"
$ws.Range("E26").Value = "  -10.38%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.18"
$ws.Range("E27").Value = "  +1.29%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "147.08"
$ws.Range("E28").Value = "  +0.89%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.036"
$ws.Range("E29").Value = "  +1.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.28"
$ws.Range("E30").Value = "  +0.71%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.746.81"
$ws.Range("E31").Value = "  -0.26%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.147"
$ws.Range("E32").Value = "  -1.64%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.008"
$ws.Range("E33").Value = "  +0.81%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9856"
$ws.Range("E34").Value = "  -4.50%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.972"
$ws.Range("E35").Value = "  -3.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08472"
$ws.Range("E36").Value = "  +0.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02548"
$ws.Range("E37").Value = "  +0.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.384"
$ws.Range("E38").Value = "  +11.07%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2315"
$ws.Range("E39").Value = "  -0.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06569"
$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.424"
$ws.Range("E41").Value = "  -2.27%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6399"
$ws.Range("E42").Value = "  +0.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.45"
$ws.Range("E43").Value = "  -3.20%  "

$ws.Range("E44").Value = "  +0.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.98"
$ws.Range("E45").Value = "  -3.47%  "

$ws.Range("E46").Value = "  +0.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5973"
$ws.Range("E47").Value = "  -0.81%  "

$ws.Range("E48").Value = "  +1.93%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.088"
$ws.Range("E49").Value = "  -2.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "125.37"
$ws.Range("E50").Value = "  +1.48%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07330"
$ws.Range("E51").Value = "  +0.66%  "
